$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.328.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.54%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.810.00'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.98%  '

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5727'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +16.30%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3882'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.58%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.18'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07600'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.98%  '

$ws.Range("E11").Value = '  +8.24%  '

$ws.Range("E12").Value = '  +0.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.15'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.71%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.252'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.50%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.807.14'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.248'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.59%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.95'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.26%  '

$ws.Range("E18").Value = '  +3.78%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06476'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.83%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.003'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.10%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.345.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.139'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.27%  '

$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.86%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.65'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.433'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +17.84%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.017.96'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.85'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.155'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +10.64%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1062'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +14.32%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.772'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.73%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.632'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2205'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.901'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +20.15%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02319'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.68%  '

$ws.Range("E38").Value = '  +6.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06116'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.99%  '

$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.039'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.32%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6391'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.10%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.163'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.56%  '

$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.379'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.06%  '

$ws.Range("E45").Value = '  +3.88%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5999'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.700'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.03'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.14%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.942'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.149'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06866'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.50%  '
